$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data in rows 2-5 with new values (100-sample batch refresh)
# Row 2
$ws.Range("A2").Value = 45053.50694444445
$ws.Range("B2").Value = 20.178
$ws.Range("C2").Value = 13.652
$ws.Range("D2").Value = 4.068
$ws.Range("E2").Value = 42.752
$ws.Range("F2").Value = 34.691
$ws.Range("G2").Value = 15.879
$ws.Range("H2").Value = 51.007
$ws.Range("I2").Value = 24.432
$ws.Range("J2").Value = 10.251
$ws.Range("K2").Value = 15.607
$ws.Range("L2").Value = 16.868
$ws.Range("M2").Value = 17.597
$ws.Range("N2").Value = 5.069
$ws.Range("O2").Value = 15.79
$ws.Range("P2").Value = 22.094
$ws.Range("Q2").Value = 13.41
$ws.Range("R2").Value = 3.498
$ws.Range("S2").Value = 2.451
$ws.Range("T2").Value = 232.996
$ws.Range("U2").Value = 43.923
$ws.Range("V2").Value = 14.575
$ws.Range("W2").Value = 29.082
$ws.Range("X2").Value = 14.996
$ws.Range("Y2").Value = 3.13
$ws.Range("Z2").Value = 25.159
$ws.Range("AA2").Value = 12.874
$ws.Range("AB2").Value = 11.654
$ws.Range("AC2").Value = 13.651
$ws.Range("AD2").Value = 17.331
$ws.Range("AE2").Value = 3.457
$ws.Range("AF2").Value = 45.214
$ws.Range("AG2").Value = 8.105
$ws.Range("AH2").Value = 18.222

# Row 3
$ws.Range("A3").Value = 45053.51388888889
$ws.Range("B3").Value = 9.609
$ws.Range("C3").Value = 6.511
$ws.Range("D3").Value = 1.581
$ws.Range("E3").Value = 20.565
$ws.Range("F3").Value = 16.566
$ws.Range("G3").Value = 7.562
$ws.Range("H3").Value = 31.653
$ws.Range("I3").Value = 11.635
$ws.Range("J3").Value = 4.908
$ws.Range("K3").Value = 7.264
$ws.Range("L3").Value = 8.229
$ws.Range("M3").Value = 8.516
$ws.Range("N3").Value = 2.418
$ws.Range("O3").Value = 7.519
$ws.Range("P3").Value = 10.522
$ws.Range("Q3").Value = 6.618
$ws.Range("R3").Value = 1.509
$ws.Range("S3").Value = 0.885
$ws.Range("T3").Value = 107.143
$ws.Range("U3").Value = 21.161
$ws.Range("V3").Value = 6.941
$ws.Range("W3").Value = 13.895
$ws.Range("X3").Value = 7.251
$ws.Range("Y3").Value = 1.576
$ws.Range("Z3").Value = 14.794
$ws.Range("AA3").Value = 6.131
$ws.Range("AB3").Value = 5.655
$ws.Range("AC3").Value = 6.607
$ws.Range("AD3").Value = 8.411
$ws.Range("AE3").Value = 1.266
$ws.Range("AF3").Value = 28.943
$ws.Range("AG3").Value = 3.8
$ws.Range("AH3").Value = 8.678

# Row 4
$ws.Range("A4").Value = 45053.52083333334
$ws.Range("B4").Value = 19.217
$ws.Range("C4").Value = 13.976
$ws.Range("D4").Value = 1.415
$ws.Range("E4").Value = 41.629
$ws.Range("F4").Value = 34.07
$ws.Range("G4").Value = 15.123
$ws.Range("H4").Value = 56.612
$ws.Range("I4").Value = 23.269
$ws.Range("J4").Value = 10.231
$ws.Range("K4").Value = 15.195
$ws.Range("L4").Value = 16.72
$ws.Range("M4").Value = 17.546
$ws.Range("N4").Value = 4.83
$ws.Range("O4").Value = 15.038
$ws.Range("P4").Value = 21.318
$ws.Range("Q4").Value = 12.771
$ws.Range("R4").Value = 1.081
$ws.Range("S4").Value = 0.907
$ws.Range("T4").Value = 221.576
$ws.Range("U4").Value = 42.004
$ws.Range("V4").Value = 13.881
$ws.Range("W4").Value = 28.128
$ws.Range("X4").Value = 14.808
$ws.Range("Y4").Value = 2.38
$ws.Range("Z4").Value = 27.556
$ws.Range("AA4").Value = 12.261
$ws.Range("AB4").Value = 10.953
$ws.Range("AC4").Value = 12.859
$ws.Range("AD4").Value = 17.425
$ws.Range("AE4").Value = 0.773
$ws.Range("AF4").Value = 51.236
$ws.Range("AG4").Value = 7.789
$ws.Range("AH4").Value = 17.354

# Row 5
$ws.Range("A5").Value = 45053.52777777778
$ws.Range("B5").Value = 21.62
$ws.Range("C5").Value = 15.9
$ws.Range("D5").Value = 1.28
$ws.Range("E5").Value = 46.92
$ws.Range("F5").Value = 38.5
$ws.Range("G5").Value = 17.01
$ws.Range("H5").Value = 65.9
$ws.Range("I5").Value = 26.18
$ws.Range("J5").Value = 11.6
$ws.Range("K5").Value = 17.22
$ws.Range("L5").Value = 18.85
$ws.Range("M5").Value = 19.83
$ws.Range("N5").Value = 5.43
$ws.Range("O5").Value = 16.92
$ws.Range("P5").Value = 24.06
$ws.Range("Q5").Value = 14.27
$ws.Range("R5").Value = 0.86
$ws.Range("S5").Value = 0.87
$ws.Range("T5").Value = 250.2
$ws.Range("U5").Value = 47.32
$ws.Range("V5").Value = 15.62
$ws.Range("W5").Value = 31.8
$ws.Range("X5").Value = 16.72
$ws.Range("Y5").Value = 2.53
$ws.Range("Z5").Value = 31.93
$ws.Range("AA5").Value = 13.79
$ws.Range("AB5").Value = 12.25
$ws.Range("AC5").Value = 14.4
$ws.Range("AD5").Value = 19.71
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 59.81
$ws.Range("AG5").Value = 8.8
$ws.Range("AH5").Value = 19.52

# Remove the now-obsolete last row (row 6)
$ws.Rows.Item(6).Delete()

# Widen columns that previously had width 7 to width 8
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
